$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($c = 1; $c -le 24; $c++) {
    Write-Host ("col " + $c + " Left=" + $ws.Cells.Item(1,$c).Left)
}
